# Shift the IGCC netting-flow timestamps forward by 22 days
# (30/31 Dec 2025 -> 21/22 Jan 2026), and keep the "Lookup" helper
# column (date + quarter-of-day number) in sync with the new dates.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item(1, 1).End(-4121).Row  # xlDown
if ($lastRow -lt 2) {
    $lastRow = 201
}

$dayShift = 22

for ($r = 2; $r -le $lastRow; $r++) {
    $tsCell = $ws.Cells.Item($r, 1)
    $quarterCell = $ws.Cells.Item($r, 4)
    $lookupCell = $ws.Cells.Item($r, 5)

    $oldTs = $tsCell.Value2
    if ($oldTs -eq $null) {
        continue
    }

    $newTs = $oldTs + $dayShift
    $tsCell.Value = $newTs

    $newDate = [datetime]::FromOADate($newTs)
    $datePart = $newDate.ToString("dd.MM.yyyy")

    $quarter = $quarterCell.Value2
    $lookupCell.Value = $datePart + $quarter
}
